# Auto-generated edit script: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.173.24"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "1.908.78"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.87"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3933"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.89"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07961"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "1.982.85"
$ws.Range("E13").Value = "  +5.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.149"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06954"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.72"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001011"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "29.213.39"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.353"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "2.173.88"
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.063"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.32"
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.864"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.003"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.45"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09423"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9249"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.363"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.260"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05850"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.166"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.016"
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02099"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5754"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1811"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.10"
$ws.Range("E44").Value = "  +2.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5413"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.224"
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07094"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.882"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.557"
$ws.Range("E49").Value = "  +6.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.10"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.075"
$ws.Range("E51").Value = "  -5.88%  "
